# Clean up garbled / mis-typed team names in the playoff data.
# (commit: "cleaned and processed data from 201819 to 201516 and started neural net code")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8  (loser of row 8): "Pittsburghsburghsburgh" -> "Pittsburgh"
$ws.Range("D8").Value = "Pittsburgh"

# Row 24 (loser of row 24): "North Carolina Asheville" -> "North Carolina-Asheville"
$ws.Range("D24").Value = "North Carolina-Asheville"

# Row 20 (loser of row 20): "University of University of University of California" -> "University of California"
$ws.Range("D20").Value = "University of California"

# Row 3  (loser of row 3): "Sothern University of University of University of California" -> "Southern California"
$ws.Range("D3").Value = "Southern California"

# Leave the cursor on the last-edited cell, matching the saved selection state.
$ws.Range("D3").Select()
